$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before D (Treatment), shifting the ROS values from D to E
$ws.Range("D1:D59").Insert(-4161)

# 2. Update header row
$ws.Range("C1").Value = "Microbiome"
$ws.Range("D1").Value = "Treatment"

# 3. Fill the new Treatment column (D) with "V" by default
$ws.Range("D2:D59").Value = "V"

# 4. Rows 45-59 (the two NAC groups) actually carry the "SPF" microbiome and a "NAC" treatment
$ws.Range("C45:C59").Value = "SPF"
$ws.Range("D45:D59").Value = "NAC"

# 5. Left-align the Genotype column (B) for all data rows + header
$ws.Range("B1:B59").HorizontalAlignment = -4131

# 6. Column widths
$ws.Columns("C").ColumnWidth = 15
$ws.Columns("D").ColumnWidth = 13.83203125

# 7. View state: scrolled down with C52 selected
$ws.Application.ActiveWindow.ScrollRow = 39
$ws.Range("C52").Select()
